$d = $word.ActiveDocument

# Locate the literal text "Adresat.name" inside the "{{ Adresat.name }}"
# placeholder (the surrounding "{{ " / " }}" runs and the gramStart/gramEnd
# proofErr markers are left untouched).
$find = $d.Content
$found = $find.Find.Execute("Adresat.name", $false, $false, $false, $false, $false,
                             $true, 1, $false, "", 0)

if ($found) {
    $target = $d.Range($find.Start, $find.End)

    # Replace "Adresat.name" with "Adresat.name" + ".text", wrapped in
    # spellStart/spellEnd proofing marks, matching what Word's spell
    # checker inserts around the newly split, out-of-dictionary token.
    $xml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
           '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
           '<pkg:xmlData>' +
           '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
           '<w:body><w:p>' +
           '<w:proofErr w:type="spellStart"/>' +
           '<w:r><w:t>Adresat.name</w:t></w:r>' +
           '<w:r><w:t>.text</w:t></w:r>' +
           '<w:proofErr w:type="spellEnd"/>' +
           '</w:p></w:body></w:document>' +
           '</pkg:xmlData></pkg:part></pkg:package>'

    $target.InsertXML($xml)
}
